# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "69.100.37"
$ws.Cells.Item(2, 5).Value = "  +1.18%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.730.58"
$ws.Cells.Item(3, 5).Value = "  +0.27%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.17%  "

# Row 5
$origStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "612.08"
$ws.Cells.Item(5, 4).Style = $origStyle
$ws.Cells.Item(5, 5).Value = "  +4.99%  "

# Row 6
$origStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "186.77"
$ws.Cells.Item(6, 4).Style = $origStyle
$ws.Cells.Item(6, 5).Value = "  +5.83%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "3.727.03"
$ws.Cells.Item(7, 5).Value = "  +0.35%  "

# Row 8
$origStyle = $ws.Cells.Item(8, 4).Style
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.640"
$ws.Cells.Item(8, 4).Style = $origStyle
$ws.Cells.Item(8, 5).Value = "  +0.66%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +0.20%  "

# Row 10
$origStyle = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.721"
$ws.Cells.Item(10, 4).Style = $origStyle
$ws.Cells.Item(10, 5).Value = "  +0.28%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -2.82%  "

# Row 12
$origStyle = $ws.Cells.Item(12, 4).Style
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "56.81"
$ws.Cells.Item(12, 4).Style = $origStyle
$ws.Cells.Item(12, 5).Value = "  +5.40%  "

# Row 13
$origStyle = $ws.Cells.Item(13, 4).Style
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000293"
$ws.Cells.Item(13, 4).Style = $origStyle
$ws.Cells.Item(13, 5).Value = "  -3.16%  "

# Row 14
$origStyle = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "10.66"
$ws.Cells.Item(14, 4).Style = $origStyle
$ws.Cells.Item(14, 5).Value = "  -0.63%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "4.326.00"
$ws.Cells.Item(15, 5).Value = "  +0.29%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "3.729.96"
$ws.Cells.Item(16, 5).Value = "  +0.06%  "

# Row 17
$origStyle = $ws.Cells.Item(17, 4).Style
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "13.09"
$ws.Cells.Item(17, 4).Style = $origStyle
$ws.Cells.Item(17, 5).Value = "  -0.25%  "

# Row 18
$origStyle = $ws.Cells.Item(18, 4).Style
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "19.38"
$ws.Cells.Item(18, 4).Style = $origStyle
$ws.Cells.Item(18, 5).Value = "  -0.53%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -0.57%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.51%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "68.942.01"
$ws.Cells.Item(21, 5).Value = "  +1.20%  "

# Row 22
$origStyle = $ws.Cells.Item(22, 4).Style
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "414.29"
$ws.Cells.Item(22, 4).Style = $origStyle
$ws.Cells.Item(22, 5).Value = "  +0.37%  "

# Row 23
$origStyle = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.63"
$ws.Cells.Item(23, 4).Style = $origStyle
$ws.Cells.Item(23, 5).Value = "  -0.35%  "

# Row 24
$origStyle = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "89.39"
$ws.Cells.Item(24, 4).Style = $origStyle
$ws.Cells.Item(24, 5).Value = "  -0.09%  "

# Row 25
$origStyle = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "3.06"
$ws.Cells.Item(25, 4).Style = $origStyle
$ws.Cells.Item(25, 5).Value = "  -1.24%  "

# Row 26
$origStyle = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "12.92"
$ws.Cells.Item(26, 4).Style = $origStyle
$ws.Cells.Item(26, 5).Value = "  -0.42%  "

# Row 27
$origStyle = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "11.03"
$ws.Cells.Item(27, 4).Style = $origStyle
$ws.Cells.Item(27, 5).Value = "  +1.77%  "

# Row 28
$origStyle = $ws.Cells.Item(28, 4).Style
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "6.06"
$ws.Cells.Item(28, 4).Style = $origStyle
$ws.Cells.Item(28, 5).Value = "  +1.99%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -0.80%  "

# Row 30
$origStyle = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "9.69"
$ws.Cells.Item(30, 4).Style = $origStyle
$ws.Cells.Item(30, 5).Value = "  +0.30%  "

# Row 31
$origStyle = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "33.38"
$ws.Cells.Item(31, 4).Style = $origStyle
$ws.Cells.Item(31, 5).Value = "  +0.38%  "

# Row 32
$origStyle = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "7.35"
$ws.Cells.Item(32, 4).Style = $origStyle
$ws.Cells.Item(32, 5).Value = "  -11.34%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -0.37%  "

# Row 34
$origStyle = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.124"
$ws.Cells.Item(34, 4).Style = $origStyle
$ws.Cells.Item(34, 5).Value = "  +4.05%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "Bittensor"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$origStyle = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "627.59"
$ws.Cells.Item(35, 4).Style = $origStyle
$ws.Cells.Item(35, 5).Value = "  +3.85%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$origStyle = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "44.94"
$ws.Cells.Item(36, 4).Style = $origStyle
$ws.Cells.Item(36, 5).Value = "  -0.67%  "

# Row 37
$origStyle = $ws.Cells.Item(37, 4).Style
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "66.33"
$ws.Cells.Item(37, 4).Style = $origStyle
$ws.Cells.Item(37, 5).Value = "  -0.08%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "0.0₃0847"
$ws.Cells.Item(38, 5).Value = "  -10.74%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +1.63%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -0.12%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.34%  "

# Row 42
$origStyle = $ws.Cells.Item(42, 4).Style
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.140"
$ws.Cells.Item(42, 4).Style = $origStyle
$ws.Cells.Item(42, 5).Value = "  +2.87%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -0.61%  "

# Row 44
$origStyle = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0446"
$ws.Cells.Item(44, 4).Style = $origStyle
$ws.Cells.Item(44, 5).Value = "  +0.27%  "

# Row 45
$origStyle = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.65"
$ws.Cells.Item(45, 4).Style = $origStyle
$ws.Cells.Item(45, 5).Value = "  +1.31%  "

# Row 46
$origStyle = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.141"
$ws.Cells.Item(46, 4).Style = $origStyle
$ws.Cells.Item(46, 5).Value = "  +3.75%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Maker"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(47, 4).Value = "2.844.97"
$ws.Cells.Item(47, 5).Value = "  +3.38%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "WEMIXToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$origStyle = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.77"
$ws.Cells.Item(48, 4).Style = $origStyle
$ws.Cells.Item(48, 5).Value = "  +5.62%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "THORChain"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$origStyle = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "9.19"
$ws.Cells.Item(49, 4).Style = $origStyle
$ws.Cells.Item(49, 5).Value = "  -4.61%  "

# Row 50
$origStyle = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.71"
$ws.Cells.Item(50, 4).Style = $origStyle
$ws.Cells.Item(50, 5).Value = "  -21.33%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "ApeXProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$origStyle = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "3.08"
$ws.Cells.Item(51, 4).Style = $origStyle
$ws.Cells.Item(51, 5).Value = "  -3.14%  "
